$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'302.81"
$ws.Range("E2").Value = "'0.99%"
$ws.Range("D3").Value = "'32.12"
$ws.Range("E3").Value = "'1.08%"
$ws.Range("D4").Value = "'4.970"
$ws.Range("E4").Value = "'-2.95%"
$ws.Range("E5").Value = "'-3.56%"
$ws.Range("D6").Value = "'2.133"
$ws.Range("E6").Value = "'-17.10%"
$ws.Range("D7").Value = "'7.861"
$ws.Range("E7").Value = "'-0.03%"
$ws.Range("D8").Value = "'3.808"
$ws.Range("E8").Value = "'-1.31%"
$ws.Range("D9").Value = "'0.9270"
$ws.Range("E9").Value = "'-0.02%"
$ws.Range("D10").Value = "'0.1755"
$ws.Range("E10").Value = "'-0.21%"
$ws.Range("D11").Value = "'0.08038"
$ws.Range("E11").Value = "'7.00%"
$ws.Range("D12").Value = "'0.08769"
$ws.Range("E12").Value = "'-3.16%"
$ws.Range("D13").Value = "'0.03161"
$ws.Range("E13").Value = "'4.78%"
$ws.Range("D14").Value = "'0.1006"
$ws.Range("E14").Value = "'0.50%"
$ws.Range("D15").Value = "'0.001510"
$ws.Range("E15").Value = "'-0.97%"
$ws.Range("D16").Value = "'0.005950"
$ws.Range("E16").Value = "'0.25%"
$ws.Range("E17").Value = "'-4.16%"
$ws.Range("D18").Value = "'2.279"
$ws.Range("E18").Value = "'-0.29%"
$ws.Range("D19").Value = "'0.3287"
$ws.Range("E19").Value = "'1.31%"
$ws.Range("D20").Value = "'0.1291"
$ws.Range("E20").Value = "'-4.17%"
$ws.Range("D21").Value = "'4.187"
$ws.Range("E21").Value = "'1.30%"
$ws.Range("D22").Value = "'0.1791"
$ws.Range("E22").Value = "'6.78%"
$ws.Range("D23").Value = "'0.04605"
$ws.Range("E23").Value = "'-0.54%"
$ws.Range("D24").Value = "'0.001237"
$ws.Range("D25").Value = "'0.004498"
$ws.Range("E25").Value = "'-1.10%"
$ws.Range("E26").Value = "'4.05%"
$ws.Range("D39").Value = "'0.01738"
$ws.Range("E39").Value = "'-2.53%"
$ws.Range("D40").Value = "'0.04803"
$ws.Range("E40").Value = "'4.61%"
$ws.Range("D41").Value = "'0.007349"
$ws.Range("E41").Value = "'6.24%"
$ws.Range("E42").Value = "'-0.76%"
$ws.Range("D43").Value = "'0.002357"
$ws.Range("E43").Value = "'10.12%"
$ws.Range("D44").Value = "'0.01112"
$ws.Range("E44").Value = "'12.89%"
$ws.Range("D45").Value = "'0.00006020"
$ws.Range("E45").Value = "'-2.64%"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("E46").Value = "'0.03%"
$ws.Range("D47").Value = "'0.003389"
$ws.Range("E47").Value = "'-59.58%"
$ws.Range("D48").Value = "'0.8234"
$ws.Range("E48").Value = "'2.17%"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("E49").Value = "'0.03%"
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("E50").Value = "'0.03%"
